$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.468.63'
$ws.Range('E2').Value = '  +1.46%  '
$ws.Range('D3').Value = '1.728.76'
$ws.Range('E3').Value = '  +2.24%  '
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.64'
$ws.Range('E5').Value = '  +2.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.0000'
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4812'
$ws.Range('E7').Value = '  +3.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2686'
$ws.Range('E8').Value = '  +2.34%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06231'
$ws.Range('E9').Value = '  +0.90%  '
$ws.Range('D10').Value = '1.730.54'
$ws.Range('E10').Value = '  +2.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07125'
$ws.Range('E11').Value = '  +1.34%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.74'
$ws.Range('E12').Value = '  +3.68%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6183'
$ws.Range('E13').Value = '  +5.91%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.554'
$ws.Range('E14').Value = '  +3.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.24'
$ws.Range('E16').Value = '  +0.27%  '
$ws.Range('D17').Value = '26.483.84'
$ws.Range('E17').Value = '  +1.52%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.0000'
$ws.Range('E18').Value = '  +0.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006963'
$ws.Range('E19').Value = '  +3.09%  '
$ws.Range('E20').Value = '  +1.79%  '
$ws.Range('D21').Value = '1.953.10'
$ws.Range('E21').Value = '  +3.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.554'
$ws.Range('E22').Value = '  +0.51%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.932'
$ws.Range('E23').Value = '  +2.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.317'
$ws.Range('E24').Value = '  +0.72%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '136.53'
$ws.Range('E25').Value = '  +1.40%  '
$ws.Range('E26').Value = '  +1.98%  '
$ws.Range('E27').Value = '  +3.50%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.405'
$ws.Range('E28').Value = '  -2.35%  '
$ws.Range('E29').Value = '  +1.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.981'
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08029'
$ws.Range('E31').Value = '  +3.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.753'
$ws.Range('E32').Value = '  +2.36%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04565'
$ws.Range('E33').Value = '  +4.47%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.617'
$ws.Range('E34').Value = '  +0.96%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6412'
$ws.Range('E35').Value = '  +4.63%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9912'
$ws.Range('E36').Value = '  +3.39%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9418'
$ws.Range('E37').Value = '  +1.18%  '
$ws.Range('E38').Value = '  +5.97%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.416'
$ws.Range('E39').Value = '  +0.91%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '107.18'
$ws.Range('E40').Value = '  -3.49%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.003'
$ws.Range('E41').Value = '  +0.65%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01501'
$ws.Range('E42').Value = '  +2.65%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.651'
$ws.Range('E43').Value = '  +11.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3923'
$ws.Range('E44').Value = '  +4.56%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.002'
$ws.Range('E45').Value = '  +12.90%  '
$ws.Range('E46').Value = '  +5.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05322'
$ws.Range('E47').Value = '  +0.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '30.98'
$ws.Range('E48').Value = '  +0.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.869'
$ws.Range('E49').Value = '  +2.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.270'
$ws.Range('E50').Value = '  +4.44%  '
$ws.Range('E51').Value = '  +2.67%  '
